# Fruta / hortaliza, semanal
# Insert a new weekly record at row 5 (pushing the existing rows 5-19 down
# to rows 6-20) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 5:19 down to 6:20, leaving a blank row 5 to fill in.
$ws.Rows("5:5").Insert()

$ws.Range("A5").Value() = 10
$ws.Range("B5").Value() = "Vega Modelo de Temuco"
$ws.Range("C5").Value() = "La Araucanía"
$ws.Range("D5").Value() = 44469
$ws.Range("E5").Value() = 9
$ws.Range("F5").Value() = 300000000
$ws.Range("G5").Value() = "Espárragos"
$ws.Range("H5").Value() = "Sin especificar"
$ws.Range("I5").Value() = "Primera"
$ws.Range("J5").Value() = 1200
$ws.Range("K5").Value() = 1800
$ws.Range("L5").Value() = 1800
$ws.Range("M5").Value() = 1800
$ws.Range("N5").Value() = "$/kilo"
$ws.Range("O5").Value() = "Región del Maule"
$ws.Range("P5").Value() = 1800
$ws.Range("Q5").Value() = 1
$ws.Range("R5").Value() = "Hortaliza"

# Match the date column's existing display format.
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
